$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.888486862182617
$ws.Range("B1").Value = 1.836326599121094
$ws.Range("C1").Value = 1.749894857406616
$ws.Range("D1").Value = 0.9783680438995361
$ws.Range("E1").Value = 0.6698675751686096
